$d = $word.ActiveDocument

# Table 4 (1-based, "permanent order" table) / Row 5 ("account_to")
# gains an explicit row height of 32 twips (1.6 pt) — i.e. a
# <w:trHeight w:val="32"/> entry in that row's <w:trPr>.
$table = $d.Tables.Item(4)
$row = $table.Rows.Item(5)
$row.Height = 1.6
